# Auto-generated cell value updates applied via Excel COM interop.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 1307.4546
$ws.Cells.Item(4, 9).Value = 985.375
$ws.Cells.Item(4, 10).Value = 2166.3333
$ws.Cells.Item(4, 11).Value = 985.375
$ws.Cells.Item(4, 12).Value = 2166.3333
$ws.Cells.Item(4, 13).Value = -871.375
$ws.Cells.Item(4, 14).Value = -2394.3333
$ws.Cells.Item(40, 8).Value = 4543
$ws.Cells.Item(40, 9).Value = 2362.75
$ws.Cells.Item(40, 10).Value = 4958.2856
$ws.Cells.Item(40, 11).Value = 2362.75
$ws.Cells.Item(40, 12).Value = 4958.2856
$ws.Cells.Item(40, 13).Value = -2187.75
$ws.Cells.Item(40, 14).Value = -5308.2856
$ws.Cells.Item(48, 8).Value = 8333.333000000001
$ws.Cells.Item(48, 10).Value = 8333.333000000001
$ws.Cells.Item(48, 12).Value = 24999.999
$ws.Cells.Item(48, 14).Value = -25583.999
$ws.Cells.Item(56, 8).Value = 8333.333000000001
$ws.Cells.Item(56, 10).Value = 8333.333000000001
$ws.Cells.Item(56, 12).Value = 24999.999
$ws.Cells.Item(56, 14).Value = -26067.999
$ws.Cells.Item(62, 8).Value = 5048.4287
$ws.Cells.Item(62, 9).Value = 1037.8
$ws.Cells.Item(62, 10).Value = 7276.5557
$ws.Cells.Item(62, 11).Value = 1037.8
$ws.Cells.Item(62, 12).Value = 7276.5557
$ws.Cells.Item(62, 13).Value = -413.8
$ws.Cells.Item(62, 14).Value = -8524.555700000001
$ws.Cells.Item(65, 8).Value = 5048.4287
$ws.Cells.Item(65, 9).Value = 1037.8
$ws.Cells.Item(65, 10).Value = 7276.5557
$ws.Cells.Item(65, 11).Value = 5189
$ws.Cells.Item(65, 12).Value = 36382.7785
$ws.Cells.Item(65, 13).Value = -2069
$ws.Cells.Item(65, 14).Value = -42622.7785
$ws.Cells.Item(69, 8).Value = 7287.3335
$ws.Cells.Item(69, 10).Value = 7287.3335
$ws.Cells.Item(69, 12).Value = 21862.0005
$ws.Cells.Item(69, 14).Value = -23610.0005
$ws.Cells.Item(72, 8).Value = 7287.3335
$ws.Cells.Item(72, 10).Value = 7287.3335
$ws.Cells.Item(72, 12).Value = 65586.0015
$ws.Cells.Item(72, 14).Value = -74322.0015
$ws.Cells.Item(80, 8).Value = 4762.8
$ws.Cells.Item(80, 9).Value = 1225
$ws.Cells.Item(80, 10).Value = 6049.273
$ws.Cells.Item(80, 11).Value = 3675
$ws.Cells.Item(80, 12).Value = 18147.819
$ws.Cells.Item(80, 13).Value = -2677
$ws.Cells.Item(80, 14).Value = -20143.819
$ws.Cells.Item(83, 8).Value = 4762.8
$ws.Cells.Item(83, 9).Value = 1225
$ws.Cells.Item(83, 10).Value = 6049.273
$ws.Cells.Item(83, 11).Value = 11025
$ws.Cells.Item(83, 12).Value = 54443.457
$ws.Cells.Item(83, 13).Value = -6033
$ws.Cells.Item(83, 14).Value = -64427.457
$ws.Cells.Item(88, 8).Value = 2415.3901
$ws.Cells.Item(88, 10).Value = 2623.5557
$ws.Cells.Item(88, 12).Value = 2623.5557
$ws.Cells.Item(88, 14).Value = -3435.5557
$ws.Cells.Item(91, 8).Value = 2415.3901
$ws.Cells.Item(91, 10).Value = 2623.5557
$ws.Cells.Item(91, 12).Value = 2623.5557
$ws.Cells.Item(91, 14).Value = -5431.5557
$ws.Cells.Item(135, 8).Value = 922.74286
$ws.Cells.Item(135, 9).Value = 507.35
$ws.Cells.Item(135, 10).Value = 1476.6
$ws.Cells.Item(135, 11).Value = 4566.150000000001
$ws.Cells.Item(135, 12).Value = 13289.4
$ws.Cells.Item(135, 13).Value = -2031.150000000001
$ws.Cells.Item(135, 14).Value = -18359.4
$ws.Cells.Item(137, 8).Value = 52851.06
$ws.Cells.Item(137, 9).Value = 67690.336
$ws.Cells.Item(137, 10).Value = 2768.5
$ws.Cells.Item(137, 11).Value = 203071.008
$ws.Cells.Item(137, 12).Value = 8305.5
$ws.Cells.Item(137, 13).Value = -200521.008
$ws.Cells.Item(137, 14).Value = -13405.5
$ws.Cells.Item(138, 8).Value = 2895.9268
$ws.Cells.Item(138, 9).Value = 1821.0869
$ws.Cells.Item(138, 11).Value = 5463.2607
$ws.Cells.Item(138, 13).Value = -323.2606999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(21, 8).Value = 8326.666999999999
$ws.Cells.Item(21, 9).Value = 5000
$ws.Cells.Item(21, 10).Value = 9990
$ws.Cells.Item(21, 11).Value = 5000
$ws.Cells.Item(21, 12).Value = 9990
$ws.Cells.Item(21, 13).Value = -4626
$ws.Cells.Item(21, 14).Value = -10738

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 4810482
$ws.Cells.Item(105, 9).Value = 5684569.5
$ws.Cells.Item(105, 11).Value = 5684569.5
$ws.Cells.Item(105, 13).Value = -5682822.5
$ws.Cells.Item(110, 8).Value = 88571.42999999999
$ws.Cells.Item(110, 10).Value = 88571.42999999999
$ws.Cells.Item(110, 12).Value = 88571.42999999999
$ws.Cells.Item(110, 14).Value = -96751.42999999999
$ws.Cells.Item(134, 8).Value = 12345.275
$ws.Cells.Item(134, 9).Value = 10510.228
$ws.Cells.Item(134, 10).Value = 18112.572
$ws.Cells.Item(134, 11).Value = 31530.684
$ws.Cells.Item(134, 12).Value = 54337.716
$ws.Cells.Item(134, 13).Value = -28995.684
$ws.Cells.Item(134, 14).Value = -59407.716
$ws.Cells.Item(140, 8).Value = 44198
$ws.Cells.Item(140, 10).Value = 43982.145
$ws.Cells.Item(140, 12).Value = 43982.145
$ws.Cells.Item(140, 14).Value = -54342.145

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 13).Value = ""
$ws.Cells.Item(31, 8).Value = 21621.79
$ws.Cells.Item(31, 9).Value = 2646.5
$ws.Cells.Item(31, 10).Value = 28612.684
$ws.Cells.Item(31, 11).Value = 2646.5
$ws.Cells.Item(31, 12).Value = 28612.684
$ws.Cells.Item(31, 13).Value = -2351.5
$ws.Cells.Item(31, 14).Value = -29202.684
$ws.Cells.Item(34, 8).Value = 21621.79
$ws.Cells.Item(34, 9).Value = 2646.5
$ws.Cells.Item(34, 10).Value = 28612.684
$ws.Cells.Item(34, 11).Value = 2646.5
$ws.Cells.Item(34, 12).Value = 28612.684
$ws.Cells.Item(34, 13).Value = -2444.5
$ws.Cells.Item(34, 14).Value = -29016.684
$ws.Cells.Item(107, 8).Value = 1613.3903
$ws.Cells.Item(107, 9).Value = 1604.8438
$ws.Cells.Item(107, 10).Value = 1643.7778
$ws.Cells.Item(107, 11).Value = 1604.8438
$ws.Cells.Item(107, 12).Value = 1643.7778
$ws.Cells.Item(107, 13).Value = 315.1561999999999
$ws.Cells.Item(107, 14).Value = -5483.7778
$ws.Cells.Item(132, 8).Value = 65535.4
$ws.Cells.Item(132, 9).Value = 57521.723
$ws.Cells.Item(132, 11).Value = 172565.169
$ws.Cells.Item(132, 13).Value = -170035.169
$ws.Cells.Item(133, 8).Value = 61484.75
$ws.Cells.Item(133, 10).Value = 63165.184
$ws.Cells.Item(133, 12).Value = 63165.184
$ws.Cells.Item(133, 14).Value = -68225.18400000001
$ws.Cells.Item(134, 8).Value = 1662.0731
$ws.Cells.Item(134, 9).Value = 1152.6285
$ws.Cells.Item(134, 10).Value = 4633.8335
$ws.Cells.Item(134, 11).Value = 3457.8855
$ws.Cells.Item(134, 12).Value = 13901.5005
$ws.Cells.Item(134, 13).Value = -922.8855000000003
$ws.Cells.Item(134, 14).Value = -18971.5005
$ws.Cells.Item(141, 8).Value = 175057.94
$ws.Cells.Item(141, 10).Value = 194932.94
$ws.Cells.Item(141, 12).Value = 194932.94
$ws.Cells.Item(141, 14).Value = -205292.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 288.6216
$ws.Cells.Item(2, 9).Value = 88.833336
$ws.Cells.Item(2, 11).Value = 533.000016
$ws.Cells.Item(2, 13).Value = -420.000016
$ws.Cells.Item(4, 8).Value = 14871579
$ws.Cells.Item(4, 9).Value = 20710414
$ws.Cells.Item(4, 11).Value = 62131242
$ws.Cells.Item(4, 13).Value = -62131130
$ws.Cells.Item(38, 8).Value = 150.72728
$ws.Cells.Item(38, 9).Value = 42.6
$ws.Cells.Item(38, 10).Value = 240.83333
$ws.Cells.Item(38, 11).Value = 127.8
$ws.Cells.Item(38, 12).Value = 722.49999
$ws.Cells.Item(38, 13).Value = 219.2
$ws.Cells.Item(38, 14).Value = -1416.49999
$ws.Cells.Item(48, 8).Value = 3834.6667
$ws.Cells.Item(48, 10).Value = 5502
$ws.Cells.Item(48, 12).Value = 16506
$ws.Cells.Item(48, 14).Value = -17006
$ws.Cells.Item(132, 8).Value = 1414.1765
$ws.Cells.Item(132, 10).Value = 1516.6666
$ws.Cells.Item(132, 12).Value = 13649.9994
$ws.Cells.Item(132, 14).Value = -18709.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 500
$ws.Cells.Item(10, 9).Value = 500
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 500
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = -331
$ws.Cells.Item(10, 14).Value = ""
$ws.Cells.Item(13, 8).Value = 5
$ws.Cells.Item(13, 9).Value = 5
$ws.Cells.Item(13, 11).Value = 5
$ws.Cells.Item(13, 13).Value = 134
$ws.Cells.Item(126, 8).Value = 4773653
$ws.Cells.Item(126, 9).Value = 3499133.5
$ws.Cells.Item(126, 11).Value = 10497400.5
$ws.Cells.Item(126, 13).Value = -10494930.5
$ws.Cells.Item(134, 8).Value = 65072
$ws.Cells.Item(134, 10).Value = 65072
$ws.Cells.Item(134, 12).Value = 195216
$ws.Cells.Item(134, 14).Value = -200286
$ws.Cells.Item(136, 8).Value = 50630
$ws.Cells.Item(136, 10).Value = 50630
$ws.Cells.Item(136, 12).Value = 151890
$ws.Cells.Item(136, 14).Value = -156990
$ws.Cells.Item(140, 8).Value = 95299.8
$ws.Cells.Item(140, 10).Value = 95299.8
$ws.Cells.Item(140, 12).Value = 95299.8
$ws.Cells.Item(140, 14).Value = -105659.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 83150
$ws.Cells.Item(22, 9).Value = 178179.6
$ws.Cells.Item(22, 10).Value = 3958.6667
$ws.Cells.Item(22, 11).Value = 178179.6
$ws.Cells.Item(22, 12).Value = 3958.6667
$ws.Cells.Item(22, 13).Value = -177884.6
$ws.Cells.Item(22, 14).Value = -4548.6667
$ws.Cells.Item(27, 8).Value = 83150
$ws.Cells.Item(27, 9).Value = 178179.6
$ws.Cells.Item(27, 10).Value = 3958.6667
$ws.Cells.Item(27, 11).Value = 178179.6
$ws.Cells.Item(27, 12).Value = 3958.6667
$ws.Cells.Item(27, 13).Value = -178072.6
$ws.Cells.Item(27, 14).Value = -4172.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 8327.25
$ws.Cells.Item(62, 9).Value = 4939
$ws.Cells.Item(62, 10).Value = 8474.565000000001
$ws.Cells.Item(62, 11).Value = 4939
$ws.Cells.Item(62, 12).Value = 8474.565000000001
$ws.Cells.Item(62, 13).Value = -4315
$ws.Cells.Item(62, 14).Value = -9722.565000000001
$ws.Cells.Item(65, 8).Value = 8327.25
$ws.Cells.Item(65, 9).Value = 4939
$ws.Cells.Item(65, 10).Value = 8474.565000000001
$ws.Cells.Item(65, 11).Value = 24695
$ws.Cells.Item(65, 12).Value = 42372.825
$ws.Cells.Item(65, 13).Value = -21575
$ws.Cells.Item(65, 14).Value = -48612.825
$ws.Cells.Item(81, 8).Value = 7250920
$ws.Cells.Item(81, 9).Value = 9807985
$ws.Cells.Item(81, 10).Value = 5903.1665
$ws.Cells.Item(81, 11).Value = 19615970
$ws.Cells.Item(81, 12).Value = 11806.333
$ws.Cells.Item(81, 13).Value = -19614909
$ws.Cells.Item(81, 14).Value = -13928.333
$ws.Cells.Item(84, 8).Value = 7250920
$ws.Cells.Item(84, 9).Value = 9807985
$ws.Cells.Item(84, 10).Value = 5903.1665
$ws.Cells.Item(84, 11).Value = 98079850
$ws.Cells.Item(84, 12).Value = 59031.665
$ws.Cells.Item(84, 13).Value = -98074546
$ws.Cells.Item(84, 14).Value = -69639.66500000001
$ws.Cells.Item(122, 8).Value = 3104.8484
$ws.Cells.Item(122, 9).Value = 1812.826
$ws.Cells.Item(122, 11).Value = 5438.478
$ws.Cells.Item(122, 13).Value = -2988.478
